$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: player played Scissors, robot played Rock, outcome is "Robot won"
$ws.Range("C2").Value = "Scissors"
$ws.Range("D2").Value = "Rock"
$ws.Range("E2").Value = "Robot won"

# Remove the now-stale rows 3 and 4 entirely (test rows)
$ws.Rows("3:4").Delete()
